# Update column F (dSF) values on Sheet1 to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F3"  = -7
    "F6"  = -5
    "F8"  = 7
    "F11" = 13
    "F13" = 6
    "F14" = -1
    "F24" = -3
    "F25" = -2
    "F30" = 3
    "F32" = -4
    "F34" = -1
    "F37" = -4
    "F40" = 5
    "F49" = -4
    "F52" = -4
    "F58" = 0
    "F59" = -3
    "F60" = 1
    "F63" = -3
    "F66" = -6
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
